# Adds the BoskampAWS weather-station rows for October 2025 (rows 33-63,
# previously blank placeholder rows) to the "September" sheet, then updates
# the sheet selection to match the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("September")

# row, day, AVG_Temperature, Max_Temperature, Min_Temperature
$rows = @(
  ,@(33, 1, 27.821527777777774, 34.4, 23.3)
  ,@(34, 2, 27.878472222222229, 33.799999999999997, 23.8)
  ,@(35, 3, 28.334751773049639, 35.299999999999997, 23.5)
  ,@(36, 4, 29.155797101449274, 36.700000000000003, 25.2)
  ,@(37, 5, 28.829861111111125, 36, 25.3)
  ,@(38, 6, 28.424999999999997, 35.4, 24.1)
  ,@(39, 7, 28.292361111111109, 34.4, 25.4)
  ,@(40, 8, 26.986805555555542, 32.9, 24.2)
  ,@(41, 9, 27.256249999999994, 33.4, 24.5)
  ,@(42, 10, 27.747916666666669, 35.5, 22.7)
  ,@(43, 11, 28.728472222222226, 36.1, 24.7)
  ,@(44, 12, 29.020833333333329, 35.5, 25.3)
  ,@(45, 13, 29.091549295774673, 36.4, 24.6)
  ,@(46, 14, 29.392142857142854, 36.4, 25.7)
  ,@(47, 15, 29.16521739130436, 35.200000000000003, 25.4)
  ,@(48, 16, 28.684615384615395, 35.299999999999997, 24.9)
  ,@(49, 17, 28.679166666666674, 35.1, 23.9)
  ,@(50, 18, 28.12222222222222, 35.799999999999997, 24.2)
  ,@(51, 19, 27.641666666666676, 35.299999999999997, 23.4)
  ,@(52, 20, 26.661594202898552, 32.6, 23.8)
  ,@(53, 21, 25.506944444444432, 29.9, 22.7)
  ,@(54, 22, 27.440579710144942, 33.9, 23)
  ,@(55, 23, 29.091666666666672, 34.799999999999997, 25.9)
  ,@(56, 24, 28.941666666666656, 35.4, 25.2)
  ,@(57, 25, 29.059027777777779, 36.1, 24.5)
  ,@(58, 26, 27.819718309859173, 34.799999999999997, 24.4)
  ,@(59, 27, 26.777536231884049, 34.299999999999997, 23.7)
  ,@(60, 28, 27.851388888888916, 33.799999999999997, 24.3)
  ,@(61, 29, 28.2326388888889, 35, 24.2)
  ,@(62, 30, 27.397916666666671, 32.700000000000003, 23.8)
  ,@(63, 31, 27.832394366197171, 33.6, 25.1)
)

foreach ($r in $rows) {
  $rowNum = $r[0]
  $day = $r[1]
  $avg = $r[2]
  $maxT = $r[3]
  $minT = $r[4]

  $ws.Cells.Item($rowNum, 1).Value = "BoskampAWS"
  $ws.Cells.Item($rowNum, 2).Value = 2025
  $ws.Cells.Item($rowNum, 3).Value = 10
  $ws.Cells.Item($rowNum, 4).Value = $day

  $ws.Cells.Item($rowNum, 5).Value = $avg
  $ws.Cells.Item($rowNum, 5).NumberFormat = "0.0"

  $ws.Cells.Item($rowNum, 6).Value = $maxT
  $ws.Cells.Item($rowNum, 6).NumberFormat = "0.0"

  $ws.Cells.Item($rowNum, 7).Value = $minT
  $ws.Cells.Item($rowNum, 7).NumberFormat = "0.0"
}

# Move the view: no more scrolled-down topLeftCell, selection now on L10.
$ws.Range("L10").Select()
